$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update "想去人数" column (F)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 85
$ws1.Range("F8").Value = 126
$ws1.Range("F9").Value = 8973
$ws1.Range("F10").Value = 822
$ws1.Range("F13").Value = 1018
$ws1.Range("F18").Value = 304
$ws1.Range("F21").Value = 1123

# Sheet "全部类型" (sheet4): same events repeated, update "想去人数" column (F)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 85
$ws4.Range("F10").Value = 126
$ws4.Range("F11").Value = 8973
$ws4.Range("F12").Value = 822
$ws4.Range("F15").Value = 1018
$ws4.Range("F20").Value = 304
$ws4.Range("F23").Value = 1123
